$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B2").Value = '108号直流'
$ws.Range("C2").Value = 46035.5384837963
$ws.Range("D2").Value = 46036.274409722224

# Row 3
$ws.Range("A3").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B3").Value = '206号直流'
$ws.Range("C3").Value = 46035.57564814815
$ws.Range("D3").Value = 46036.274409722224

# Row 4
$ws.Range("A4").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B4").Value = '103号直流'
$ws.Range("C4").Value = 46035.57846064815
$ws.Range("D4").Value = 46036.274409722224

# Row 5
$ws.Range("A5").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B5").Value = '106号直流'
$ws.Range("C5").Value = 46035.6333912037
$ws.Range("D5").Value = 46036.274409722224

# Row 6
$ws.Range("A6").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B6").Value = '204号直流'
$ws.Range("C6").Value = 46035.67798611111
$ws.Range("D6").Value = 46036.274409722224

# Row 7
$ws.Range("A7").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B7").Value = '212号直流'
$ws.Range("C7").Value = 46035.68784722222
$ws.Range("D7").Value = 46036.274409722224

# Row 8
$ws.Range("A8").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B8").Value = '105号直流'
$ws.Range("C8").Value = 46035.709756944445
$ws.Range("D8").Value = 46036.274409722224

# Row 9
$ws.Range("A9").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B9").Value = '309号直流'
$ws.Range("C9").Value = 46035.74162037037
$ws.Range("D9").Value = 46036.274409722224

# Row 10
$ws.Range("A10").Value = '飞狐四方坪西区充电站'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = '9176699400500405'
$ws.Range("A10").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C10").Value = 46035.42306712963
$ws.Range("D10").Value = 46036.32145833333

# Row 11
$ws.Range("A11").Value = '飞狐四方坪西区充电站'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = '9176699400500205'
$ws.Range("A11").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C11").Value = 46035.53597222222
$ws.Range("D11").Value = 46036.32145833333

# Row 12
$ws.Range("A12").Value = '飞狐四方坪西区充电站'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = '9176699400500102'
$ws.Range("A12").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C12").Value = 46035.538298611114
$ws.Range("D12").Value = 46036.32145833333

# Row 13
$ws.Range("A13").Value = '飞狐四方坪西区充电站'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = '9176699400500601'
$ws.Range("A13").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C13").Value = 46035.53974537037
$ws.Range("D13").Value = 46036.32145833333

# Row 14
$ws.Range("A14").Value = '飞狐四方坪南区充电站'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = '9176699368200102'
$ws.Range("A14").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C14").Value = 46035.54384259259
$ws.Range("D14").Value = 46036.32145833333

# Row 15
$ws.Range("A15").Value = '飞狐四方坪西区充电站'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = '9176699400501205'
$ws.Range("A15").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C15").Value = 46035.55332175926
$ws.Range("D15").Value = 46036.32145833333

# Row 16
$ws.Range("A16").Value = '飞狐四方坪西区充电站'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = '9176699400500403'
$ws.Range("A16").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C16").Value = 46035.55542824074
$ws.Range("D16").Value = 46036.32145833333

# Row 17
$ws.Range("A17").Value = '飞狐四方坪东区充电站'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = '9176699425700301'
$ws.Range("A17").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C17").Value = 46035.559166666666
$ws.Range("D17").Value = 46036.32145833333

# Row 18
$ws.Range("A18").Value = '飞狐四方坪西区充电站'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = '9176699400500604'
$ws.Range("A18").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C18").Value = 46035.56146990741
$ws.Range("D18").Value = 46036.32145833333

# Row 19
$ws.Range("A19").Value = '飞狐四方坪西区充电站'
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = '9176699400501304'
$ws.Range("A19").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C19").Value = 46035.5653125
$ws.Range("D19").Value = 46036.32145833333

# Row 20
$ws.Range("A20").Value = '飞狐四方坪西区充电站'
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = '9176699400500201'
$ws.Range("A20").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C20").Value = 46035.56642361111
$ws.Range("D20").Value = 46036.32145833333

# Row 21
$ws.Range("A21").Value = '飞狐四方坪西区充电站'
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = '9176699400501303'
$ws.Range("A21").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C21").Value = 46035.5690625
$ws.Range("D21").Value = 46036.32145833333

# Row 22
$ws.Range("A22").Value = '飞狐四方坪东区充电站'
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = '9176699416300203'
$ws.Range("A22").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C22").Value = 46035.57331018519
$ws.Range("D22").Value = 46036.32145833333

# Row 23
$ws.Range("A23").Value = '飞狐四方坪西区充电站'
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = '9176699400501002'
$ws.Range("A23").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C23").Value = 46035.575150462966
$ws.Range("D23").Value = 46036.32145833333

# Row 24
$ws.Range("A24").Value = '飞狐四方坪南区充电站'
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = '9176699368200101'
$ws.Range("A24").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C24").Value = 46035.575740740744
$ws.Range("D24").Value = 46036.32145833333

# Row 25
$ws.Range("A25").Value = '飞狐四方坪西区充电站'
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = '9176699400501203'
$ws.Range("A25").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C25").Value = 46035.57679398148
$ws.Range("D25").Value = 46036.32145833333

# Row 26
$ws.Range("A26").Value = '飞狐四方坪西区充电站'
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = '9176699400501101'
$ws.Range("A26").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C26").Value = 46035.57934027778
$ws.Range("D26").Value = 46036.32145833333

# Row 27
$ws.Range("A27").Value = '飞狐四方坪南区充电站'
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = '9176699368200406'
$ws.Range("A27").Copy()
$ws.Range("B27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C27").Value = 46035.58269675926
$ws.Range("D27").Value = 46036.32145833333

# Row 28
$ws.Range("A28").Value = '飞狐四方坪南区充电站'
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = '9176699368200103'
$ws.Range("A28").Copy()
$ws.Range("B28").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C28").Value = 46035.583032407405
$ws.Range("D28").Value = 46036.32145833333

# Row 29
$ws.Range("A29").Value = '飞狐四方坪东区充电站'
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = '9176699442100702'
$ws.Range("A29").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C29").Value = 46035.644166666665
$ws.Range("D29").Value = 46036.32145833333

# Row 30
$ws.Range("A30").Value = '飞狐四方坪东区充电站'
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = '9176699442100302'
$ws.Range("A30").Copy()
$ws.Range("B30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C30").Value = 46035.66125
$ws.Range("D30").Value = 46036.32145833333

# Row 31
$ws.Range("A31").Value = '飞狐四方坪南区充电站'
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = '9176699368200203'
$ws.Range("A31").Copy()
$ws.Range("B31").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C31").Value = 46035.68814814815
$ws.Range("D31").Value = 46036.32145833333

# Row 32
$ws.Range("A32").Value = '飞狐四方坪西区充电站'
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = '9176699400500305'
$ws.Range("A32").Copy()
$ws.Range("B32").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C32").Value = 46035.718298611115
$ws.Range("D32").Value = 46036.32145833333

# Row 33
$ws.Range("A33").Value = '飞狐四方坪东区充电站'
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = '9176699442100901'
$ws.Range("A33").Copy()
$ws.Range("B33").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C33").Value = 46035.80175925926
$ws.Range("D33").Value = 46036.32145833333

# Clear rows 34-36 (A:D) entirely
$ws.Range("A34:D36").ClearContents()

# Update selection to match target view state
$ws.Range("E15").Select()
